$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5 use "Ano " prefix for B1:E1 (years 2015, 2030, 2040, 2050)
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4 uses "Intervalo " prefix for B1:E1 (2015, 2015-2030, 2031-2040, 2041-2050)
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo 2015"
$ws.Range("C1").Value = "Intervalo 2015-2030"
$ws.Range("D1").Value = "Intervalo 2031-2040"
$ws.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 6 only has B1 (2015) -> "Ano 2015"
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano 2015"
